$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "○" check mark for IntroScene row (row 4)
$ws.Range("B4").Value = "○"

# Add a new row for AtHomeScene with its own check mark
$ws.Range("A5").Value = "AtHomeScene"
$ws.Range("B5").Value = "○"

# Leave the selection where the author last clicked
$ws.Range("G9").Select()
